# Clean up code and fix output
# Add a new "Yearly demand" worksheet at the end of the workbook containing
# the net yearly demand profile (3 representative days x 24 hours).

$wb = $excel.ActiveWorkbook

# Use the first sheet's header-style cell (B1) as the formatting template:
# bold font, thin border, centered horizontal / top vertical alignment.
$styleTemplate = $wb.Worksheets.Item(1).Range("B1")

# Add the new worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "Yearly demand"

# --- Header row (B1:Y1) = 0..23 (hour index) ---
for ($i = 0; $i -le 23; $i++) {
    $newSheet.Cells.Item(1, 2 + $i).Value = $i
}

# --- Column A (A2=0, A3=1, A4=2) day index ---
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(4, 1).Value = 2

# --- Data rows ---
$row2 = @(-32.5, -19.5, -13, -13, -13, 142.5, 291.5, 327, 388.5, 502, 596, 670.5, 745, 651, 576.5, 502, 320.5, 139, 32, -117, -97.5, -78, -52, -39)
$row3 = @(-32.5, -19.5, -13, 0, 0, -19.5, 0, 324, 486, 648, 729, 751.5, 583, 567, 333.5, 340, 243, 57.99999999999999, -130, 0, 0, -78, 0, -39)
$row4 = @(-32.5, -19.5, 0, 0, 0, -19.5, 0, 0, 81, 324, 567, 589.5, 648, 567, 324, 162, 81, 0, -130, 0, 0, 0, 0, -39)

for ($i = 0; $i -le 23; $i++) {
    $newSheet.Cells.Item(2, 2 + $i).Value = $row2[$i]
    $newSheet.Cells.Item(3, 2 + $i).Value = $row3[$i]
    $newSheet.Cells.Item(4, 2 + $i).Value = $row4[$i]
}

# --- Apply the same header/label style used on the other sheets (style id 1:
#     bold, bordered, centered-top) by copying formats from the template cell ---
$styleTemplate.Copy()
$newSheet.Range("B1:Y1").PasteSpecial(-4122)   # xlPasteFormats
$newSheet.Range("A2:A4").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# Reset selection/active cell to A1 to match the other sheets' saved view state.
$newSheet.Range("A1").Select()
